$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# --- Row 21: fill in the MELD dataset's data (A21 "MELD" already present) ---
$ws.Range("B21").Value = "acted"
$ws.Range("C21").Value = 3086
$ws.Range("D21").Value = 4186
$ws.Range("E21").Value = 6430
$ws.Range("F21").Value = "English"
$ws.Range("H21").Value = 356
$ws.Range("I21").Value = "Joy, Sadness, Fear, Anger, Surprise, Disgust, Neutral"

# J21: rich text "Friends" (cell default formatting) + " data" (explicit, non-italic)
$ws.Range("J21").Value = "Friends data"
$ws.Range("J21").Font.Italic = $true
$ws.Range("J21").Characters(8, 5).Font.Italic = $false

# --- Row 22: start of the next dataset entry ("oreau2", French) ---
$ws.Range("A22").Value = "oreau2"
$ws.Range("F22").Value = "French"

# --- Sheet view / print setup touch-ups ---
$ws.PageSetup.Orientation = 1
$ws.Range("A22").Select() | Out-Null
